$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A56").Value = 41934
$ws.Range("A56").NumberFormat = "m/d/yyyy"
$ws.Range("B56").Value = "Forum"
$ws.Range("C56").Value = "Stephen.Frank@nrel.gov"
$ws.Range("D56").Value = "Request to make it easier to transfer sam inputs to SDK -- enhancements to inputs browser"
$ws.Range("F56").Value = 41934
$ws.Range("F56").NumberFormat = "m/d/yyyy"

$ws.Rows.Item(56).RowHeight = 30

$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("A57").Select()
